# Updated cryptos list with GitHub Actions
# Prices in column D are stored as text (they use "." as a thousands
# separator, e.g. "24.386.35"), so every assignment is prefixed with an
# apostrophe to force Excel to keep them as text instead of auto-coercing
# plain-decimal-looking values (e.g. "1.002") into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "'24.386.35"
$ws.Cells.Item(2, 5).Value = "  -1.73%  "

# row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "'1.653.66"
$ws.Cells.Item(3, 5).Value = "  -2.77%  "

# row 4 - TetherUSD
$ws.Cells.Item(4, 4).Value = "'1.002"
$ws.Cells.Item(4, 5).Value = "  +0.19%  "

# row 5 - BNB
$ws.Cells.Item(5, 4).Value = "'311.33"
$ws.Cells.Item(5, 5).Value = "  -1.10%  "

# row 6 - USDC
$ws.Cells.Item(6, 5).Value = "  +0.05%  "

# row 7 - XRP
$ws.Cells.Item(7, 4).Value = "'0.3912"
$ws.Cells.Item(7, 5).Value = "  -1.87%  "

# row 8 - Cardano
$ws.Cells.Item(8, 4).Value = "'0.3914"
$ws.Cells.Item(8, 5).Value = "  -3.01%  "

# row 9 - BinanceUSD
$ws.Cells.Item(9, 5).Value = "  +0.28%  "

# row 10 - Polygon
$ws.Cells.Item(10, 5).Value = "  -6.11%  "

# row 11 - OKB
$ws.Cells.Item(11, 4).Value = "'49.88"
$ws.Cells.Item(11, 5).Value = "  -6.86%  "

# row 12 - Dogecoin
$ws.Cells.Item(12, 4).Value = "'0.08541"
$ws.Cells.Item(12, 5).Value = "  -2.91%  "

# row 13 - Solana
$ws.Cells.Item(13, 4).Value = "'24.81"
$ws.Cells.Item(13, 5).Value = "  -4.87%  "

# row 14 - Polkadot
$ws.Cells.Item(14, 5).Value = "  -4.36%  "

# row 15 - ShibaInu
$ws.Cells.Item(15, 5).Value = "  -3.17%  "

# row 16 - Chainlink
$ws.Cells.Item(16, 5).Value = "  -4.77%  "

# row 17 - WrappedEther
$ws.Cells.Item(17, 4).Value = "'1.661.24"
$ws.Cells.Item(17, 5).Value = "  -4.17%  "

# row 18 - Litecoin
$ws.Cells.Item(18, 4).Value = "'93.46"
$ws.Cells.Item(18, 5).Value = "  -2.26%  "

# row 19 - TRON
$ws.Cells.Item(19, 4).Value = "'0.06944"
$ws.Cells.Item(19, 5).Value = "  -3.26%  "

# row 20 - Avalanche
$ws.Cells.Item(20, 4).Value = "'20.87"
$ws.Cells.Item(20, 5).Value = "  -0.11%  "

# row 21 - Uniswap
$ws.Cells.Item(21, 4).Value = "'7.011"
$ws.Cells.Item(21, 5).Value = "  -4.62%  "

# row 22 - Dai
$ws.Cells.Item(22, 5).Value = "  +0.06%  "

# row 23 - Cosmos
$ws.Cells.Item(23, 5).Value = "  -4.02%  "

# row 24 - WrappedBTC
$ws.Cells.Item(24, 4).Value = "'24.394.43"
$ws.Cells.Item(24, 5).Value = "  -1.56%  "

# row 25 - Toncoin
$ws.Cells.Item(25, 4).Value = "'2.340"
$ws.Cells.Item(25, 5).Value = "  -1.16%  "

# row 26 - LidoDAOToken
$ws.Cells.Item(26, 4).Value = "'2.779"
$ws.Cells.Item(26, 5).Value = "  -4.89%  "

# row 27 - EthereumClassic
$ws.Cells.Item(27, 4).Value = "'22.69"
$ws.Cells.Item(27, 5).Value = "  -2.01%  "

# row 28 - Monero
$ws.Cells.Item(28, 4).Value = "'159.12"
$ws.Cells.Item(28, 5).Value = "  -1.63%  "

# row 29 - HuobiToken
$ws.Cells.Item(29, 4).Value = "'5.682"
$ws.Cells.Item(29, 5).Value = "  -7.64%  "

# row 30 - BitcoinCash
$ws.Cells.Item(30, 4).Value = "'145.18"
$ws.Cells.Item(30, 5).Value = "  +0.82%  "

# row 31 - Filecoin
$ws.Cells.Item(31, 4).Value = "'8.167"
$ws.Cells.Item(31, 5).Value = "  -3.07%  "

# row 32 - WEMIXTOKEN
$ws.Cells.Item(32, 4).Value = "'2.586"
$ws.Cells.Item(32, 5).Value = "  +9.50%  "

# row 33 - WrappedliquidstakedEther2.0
$ws.Cells.Item(33, 4).Value = "'1.840.69"
$ws.Cells.Item(33, 5).Value = "  -1.07%  "

# row 34 - ImmutableX
$ws.Cells.Item(34, 4).Value = "'1.013"
$ws.Cells.Item(34, 5).Value = "  -2.23%  "

# row 35 - now Hedera (was VeChain)
$ws.Cells.Item(35, 2).Value = "Hedera"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(35, 4).Value = "'0.08158"
$ws.Cells.Item(35, 5).Value = "  -5.66%  "

# row 36 - now VeChain (was Hedera)
$ws.Cells.Item(36, 2).Value = "VeChain"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(36, 4).Value = "'0.03020"
$ws.Cells.Item(36, 5).Value = "  -4.82%  "

# row 37 - InternetComputer(DFINITY)
$ws.Cells.Item(37, 4).Value = "'6.854"
$ws.Cells.Item(37, 5).Value = "  -6.71%  "

# row 38 - Algorand
$ws.Cells.Item(38, 4).Value = "'0.2766"
$ws.Cells.Item(38, 5).Value = "  -2.43%  "

# row 39 - Stellar
$ws.Cells.Item(39, 4).Value = "'0.09437"
$ws.Cells.Item(39, 5).Value = "  -0.31%  "

# row 41 - TrustWalletToken
$ws.Cells.Item(41, 4).Value = "'1.488"
$ws.Cells.Item(41, 5).Value = "  +0.72%  "

# row 42 - TheSandbox
$ws.Cells.Item(42, 4).Value = "'0.7807"
$ws.Cells.Item(42, 5).Value = "  -6.05%  "

# row 43 - Aptos
$ws.Cells.Item(43, 4).Value = "'13.43"

# row 44 - EnergySwap
$ws.Cells.Item(44, 4).Value = "'16.35"
$ws.Cells.Item(44, 5).Value = "  -7.42%  "

# row 45 - NEARProtocol
$ws.Cells.Item(45, 4).Value = "'2.554"
$ws.Cells.Item(45, 5).Value = "  -5.59%  "

# row 46 - Decentraland
$ws.Cells.Item(46, 4).Value = "'0.7030"
$ws.Cells.Item(46, 5).Value = "  -5.42%  "

# row 47 - PancakeSwap
$ws.Cells.Item(47, 4).Value = "'4.149"
$ws.Cells.Item(47, 5).Value = "  -1.51%  "

# row 48 - Cronos
$ws.Cells.Item(48, 4).Value = "'0.08616"
$ws.Cells.Item(48, 5).Value = "  +2.90%  "

# row 49 - Frax
$ws.Cells.Item(49, 4).Value = "'1.001"
$ws.Cells.Item(49, 5).Value = "  +0.12%  "

# row 50 - Flow
$ws.Cells.Item(50, 4).Value = "'1.309"
$ws.Cells.Item(50, 5).Value = "  -5.55%  "

# row 51 - Quant
$ws.Cells.Item(51, 4).Value = "'136.71"
$ws.Cells.Item(51, 5).Value = "  -2.22%  "
